# Append a new record (row 60) to the "Sheet1" worksheet (the tab backed by
# xl/worksheets/sheet2.xml) in Combined.xlsx.
#
# All of the existing data rows store every value as text (t="str"/shared
# string), even when the text looks numeric ("14", "1452", ...). Assigning
# such strings straight to .Value would make Excel auto-coerce them into
# numbers, so we briefly force Text number-format before writing, then clear
# the formatting again so the new cells end up back on the sheet's default
# (unstyled) cell style -- exactly like every other cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rowNum = 60
$targetRange = $ws.Range("A60:AJ60")

# Force text storage so numeric-looking strings aren't reinterpreted as numbers.
$targetRange.NumberFormat = "@"

$values = @(
  "DIFF",
  "xyz",
  "KLOPD3412",
  "-",
  "3",
  "14",
  "12",
  "14",
  "30x8x2.5",
  "424",
  "-",
  "4",
  "Black",
  "Alloy steel",
  "choice of screws or glue + Screw cap",
  "-",
  "Danpoo",
  "14.3",
  "1",
  "1452",
  "100",
  "14,854.08",
  "5236",
  "4",
  "133,172",
  "CN",
  "1425",
  "452",
  "474",
  "5",
  "MEISHIMEIQI US",
  "Home & Kitchen",
  "Large and Bulky",
  "FBA",
  "2026-01-01",
  "Danpoo"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($rowNum, $col).Value = $values[$i]
}

# Drop the temporary Text number-format so the new cells fall back to the
# workbook's default (unstyled) cell format, matching the rest of the sheet.
$targetRange.ClearFormats()
